$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row 50 data (quincena-pasada adjustment for Johana Quimbay)
# Force column A to stay textual ("2024-10-05") instead of being
# auto-converted into a date serial number, matching the existing
# text-formatted date cells in rows 45-49.
$ws.Range("A50").NumberFormat = "@"
$ws.Range("A50").Value = "2024-10-05"
$ws.Range("A50").Style = "Normal"

$ws.Range("C50").Value = "Descuento - Producto - Ajuste Quincena Pasada"
$ws.Range("D50").Value = "Johana Quimbay"
$ws.Range("F50").Value = -27000
